$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-76 down to 27-77.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly price entry.
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 45044
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 100112042
$ws.Range("G26").Value = "Locoto"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 4400
$ws.Range("L26").Value = 4400
$ws.Range("M26").Value = 4400
$ws.Range("N26").Value = "$/kilo"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 4400
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"
